$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new quarterly columns before column D. This shifts the existing
# D:K data right to F:M, matching the source workbook's new layout
# (two additional quarters: 2018-12-31 and 2018-09-23 prepended).
# ---------------------------------------------------------------------------
$ws.Range("D:E").Insert()

# Propagate cell formatting (number format / font / alignment) from the
# shifted-right column F (which now carries the original column D's
# formatting) onto the two newly inserted columns D:E, for every row used
# by the financial statement tables.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Populate the two new columns with the latest quarter data.
# Each entry is: row, column D value, column E value ($null = leave blank).
# ---------------------------------------------------------------------------
$newData = @(
    ,@(7, 43465, 43366)
    ,@(8, 249800, 663700)
    ,@(9, 19800, 53900)
    ,@(10, 230000, 609800)
    ,@(11, $null, $null)
    ,@(12, "NA", "NA")
    ,@(13, 0, 0)
    ,@(14, 2200, 3200)
    ,@(15, 23400, 74400)
    ,@(16, $null, $null)
    ,@(17, 210400, 405100)
    ,@(18, 39400, 258600)
    ,@(19, $null, $null)
    ,@(20, -37300, 15000)
    ,@(21, 25500, 347900)
    ,@(22, 23100, 21500)
    ,@(23, -21100, 252100)
    ,@(24, 11300, 38800)
    ,@(25, 0, 0)
    ,@(26, -32400, 213300)
    ,@(27, -32400, 213300)
    ,@(28, 0, 0)
    ,@(29, 9900, "NA")
    ,@(30, 0, 0)
    ,@(31, 0, 0)
    ,@(32, 37300, -15000)
    ,@(33, -22500, 213300)
    ,@(34, 0, 0)
    ,@(35, -22500, 213300)
    ,@(38, 43465, 43366)
    ,@(39, $null, $null)
    ,@(40, $null, $null)
    ,@(41, 105300, 190800)
    ,@(42, 0, 0)
    ,@(43, 51500, 58400)
    ,@(44, 30800, 36500)
    ,@(45, 12600, 21900)
    ,@(46, 200200, 307600)
    ,@(47, 0, 0)
    ,@(48, 1599400, 1588700)
    ,@(49, 215100, 219100)
    ,@(50, 0, 0)
    ,@(51, 0, 0)
    ,@(52, 9400, 13500)
    ,@(53, 0, 0)
    ,@(54, 2024200, 2128900)
    ,@(55, $null, $null)
    ,@(56, $null, $null)
    ,@(57, 23300, 33000)
    ,@(58, 5600, 3800)
    ,@(59, 205800, 251300)
    ,@(60, 234700, 288100)
    ,@(61, 1657600, 1658300)
    ,@(62, 99500, 90900)
    ,@(63, 0, 0)
    ,@(64, 0, 0)
    ,@(65, 0, 0)
    ,@(66, 1991800, 2037300)
    ,@(67, $null, $null)
    ,@(68, 0, 0)
    ,@(69, 0, 0)
    ,@(70, 0, 0)
    ,@(71, 0, 0)
    ,@(72, 0, 0)
    ,@(73, 0, 0)
    ,@(74, 0, 0)
    ,@(75, 0, 0)
    ,@(76, 32400, 91700)
    ,@(77, 0, 0)
    ,@(80, 43465, 43366)
    ,@(81, -22500, 213300)
    ,@(82, $null, $null)
    ,@(83, 23400, 74400)
    ,@(84, 0, 0)
    ,@(85, 0, 0)
    ,@(86, 0, 0)
    ,@(87, 0, 0)
    ,@(88, 0, 0)
    ,@(89, 16300, 248700)
    ,@(90, $null, $null)
    ,@(91, -44100, -45100)
    ,@(92, 0, 0)
    ,@(93, 0, 0)
    ,@(94, -44000, -45100)
    ,@(95, $null, $null)
    ,@(96, -52300, -50300)
    ,@(97, 0, 0)
    ,@(98, 0, 0)
    ,@(99, 0, 0)
    ,@(100, -53300, -75300)
    ,@(101, -4400, 2400)
    ,@(102, -85400, 130600)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    if ($null -ne $dVal) {
        $ws.Cells.Item($r, 4).Value = $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($r, 5).Value = $eVal
    }
}

# ---------------------------------------------------------------------------
# Correct two pre-existing data entries in the "Capital Expenditures" row
# (now row 91 after the shift) that were wrong in the source data.
# ---------------------------------------------------------------------------
$ws.Range("I91").Value = -28700
$ws.Range("J91").Value = -75200
